$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Terminal Hortofrutícola Agro Chillán - Mango,
# Primera, 2021-11-05) was reported and belongs right after the existing
# row 56 (2021-10-02) entry, so insert a fresh row at 57 and push the
# following historical rows (old 57-65) down to 58-66.
$ws.Rows("57:57").Insert()

$ws.Cells.Item(57, 1).Value = 7
$ws.Cells.Item(57, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(57, 3).Value = "Ñuble"
$ws.Cells.Item(57, 4).Value = 44505
$ws.Cells.Item(57, 5).Value = 16
$ws.Cells.Item(57, 6).Value = "Fruta"
$ws.Cells.Item(57, 7).Value = 100108
$ws.Cells.Item(57, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(57, 9).Value = 100108002
$ws.Cells.Item(57, 10).Value = "Mango"
$ws.Cells.Item(57, 11).Value = "Sin especificar"
$ws.Cells.Item(57, 12).Value = "Primera"
$ws.Cells.Item(57, 13).Value = 60
$ws.Cells.Item(57, 14).Value = 7500
$ws.Cells.Item(57, 15).Value = 8000
$ws.Cells.Item(57, 16).Value = 7750
$ws.Cells.Item(57, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(57, 18).Value = "Perú"
$ws.Cells.Item(57, 19).Value = 1938
$ws.Cells.Item(57, 20).Value = 4
